$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (weekly update: new record added at top of
# the data, pushing existing rows 4-26 down to 5-27).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest week's record.
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 45092
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104003
$ws.Range("J4").Value = "Membrillo"
$ws.Range("K4").Value = "Champion"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 110
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10455
$ws.Range("Q4").Value = "$/bandeja 18 kilos granel"
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 581
$ws.Range("T4").Value = 18
